$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.995.67'
$ws.Range("E2").Value = '  -0.22%  '
$ws.Range("D3").Value = '2.227.13'
$ws.Range("E3").Value = '  -0.36%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '''252.42'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.66%  '
$ws.Range("D6").Value = '''0.632'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.78%  '
$ws.Range("D7").Value = '''71.98'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +5.24%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("E9").Value = '  +10.15%  '
$ws.Range("D10").Value = '''40.53'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +13.67%  '
$ws.Range("D11").Value = '''0.0973'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.47%  '
$ws.Range("D12").Value = '''58.33'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.70%  '
$ws.Range("D13").Value = '''7.30'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +9.73%  '
$ws.Range("D14").Value = '''0.105'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.72%  '
$ws.Range("D15").Value = '2.559.40'
$ws.Range("E15").Value = '  -0.45%  '
$ws.Range("E16").Value = '  +1.59%  '
$ws.Range("D17").Value = '''0.892'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +4.24%  '
$ws.Range("D18").Value = '2.225.01'
$ws.Range("E18").Value = '  -0.56%  '
$ws.Range("D19").Value = '41.939.65'
$ws.Range("E19").Value = '  -0.22%  '
$ws.Range("D20").Value = '0.0₃0971'
$ws.Range("E20").Value = '  +1.47%  '
$ws.Range("E21").Value = '  +1.90%  '
$ws.Range("D22").Value = '''72.94'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.08%  '
$ws.Range("D23").Value = '''235.95'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.73%  '
$ws.Range("D24").Value = '''2.10'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.78%  '
$ws.Range("D25").Value = '''4.11'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +13.26%  '
$ws.Range("D26").Value = '''12.04'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +21.81%  '
$ws.Range("E27").Value = '  +0.11%  '
$ws.Range("D28").Value = '''2.54'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.45%  '
$ws.Range("E29").Value = '  -1.25%  '
$ws.Range("D30").Value = '''170.99'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.40%  '
$ws.Range("D31").Value = '''20.92'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.84%  '
$ws.Range("D32").Value = '''0.124'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.28%  '
$ws.Range("D33").Value = '''5.60'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +7.28%  '
$ws.Range("E34").Value = '  -0.52%  '
$ws.Range("D35").Value = '''0.0746'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +5.01%  '
$ws.Range("E36").Value = '  +0.98%  '
$ws.Range("D37").Value = '''26.34'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +17.64%  '
$ws.Range("D38").Value = '''4.14'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +10.70%  '
$ws.Range("D39").Value = '''0.0309'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +7.67%  '
$ws.Range("D40").Value = '''2.29'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.29%  '
$ws.Range("D41").Value = '''5.96'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.01%  '
$ws.Range("D42").Value = '''12.51'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +26.92%  '
$ws.Range("D43").Value = '''66.03'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.48%  '
$ws.Range("D44").Value = '''0.207'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +9.66%  '
$ws.Range("D45").Value = '''4.87'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.59%  '
$ws.Range("D46").Value = '''8.89'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.67%  '
$ws.Range("E47").Value = '  +1.01%  '
$ws.Range("D48").Value = '''4.67'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.43%  '
$ws.Range("E49").Value = '  -0.17%  '
$ws.Range("D50").Value = '''1.17'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +7.33%  '
$ws.Range("E51").Value = '  +6.14%  '
